# Auto-generated edit script: updates crypto price/volume cells
# to match the latest scrape, per the GitHub Actions commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "63.481.02"
$ws.Range("E2").Value = "  -0.15%  "
$ws.Range("D3").Value = "2.616.59"
$ws.Range("E3").Value = "  -1.13%  "
$ws.Range("E4").Value = "  +0.17%  "
$ws.Range("D5").Value = "'594.35"
$ws.Range("E5").Value = "  -1.47%  "
$ws.Range("D6").Value = "'149.81"
$ws.Range("E6").Value = "  +2.20%  "
$ws.Range("E7").Value = "  +0.09%  "
$ws.Range("D8").Value = "'0.587"
$ws.Range("E8").Value = "  -0.49%  "
$ws.Range("D9").Value = "'0.108"
$ws.Range("E9").Value = "  +0.26%  "
$ws.Range("D10").Value = "'5.66"
$ws.Range("E10").Value = "  +1.39%  "
$ws.Range("D11").Value = "'0.383"
$ws.Range("E11").Value = "  +3.67%  "
$ws.Range("E12").Value = "  -1.31%  "
$ws.Range("D13").Value = "'27.60"
$ws.Range("E13").Value = "  +0.29%  "
$ws.Range("D14").Value = "3.089.80"
$ws.Range("E14").Value = "  -0.98%  "
$ws.Range("D15").Value = "63.345.43"
$ws.Range("E15").Value = "  -0.16%  "
$ws.Range("E16").Value = "  +1.80%  "
$ws.Range("D17").Value = "2.613.65"
$ws.Range("E17").Value = "  -0.56%  "
$ws.Range("E18").Value = "  +6.91%  "
$ws.Range("D19").Value = "'4.62"
$ws.Range("E19").Value = "  +1.53%  "
$ws.Range("D20").Value = "'346.61"
$ws.Range("E20").Value = "  +1.40%  "
$ws.Range("D21").Value = "'6.84"
$ws.Range("E21").Value = "  -1.33%  "
$ws.Range("E22").Value = "  -0.21%  "
$ws.Range("D23").Value = "'5.69"
$ws.Range("E23").Value = "  +1.92%  "
$ws.Range("D24").Value = "'66.23"
$ws.Range("E24").Value = "  -0.66%  "
$ws.Range("E25").Value = "  +12.15%  "
$ws.Range("B26").Value = "Fetch.AI"
$ws.Range("C26").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D26").Value = "'1.67"
$ws.Range("E26").Value = "  -1.12%  "
$ws.Range("B27").Value = "InternetComputer(DFINITY)"
$ws.Range("C27").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D27").Value = "'9.14"
$ws.Range("E27").Value = "  +0.62%  "
$ws.Range("D28").Value = "'566.39"
$ws.Range("E28").Value = "  +0.20%  "
$ws.Range("D29").Value = "'8.20"
$ws.Range("E29").Value = "  +2.73%  "
$ws.Range("E30").Value = "  -0.10%  "
$ws.Range("E31").Value = "  +0.24%  "
$ws.Range("D33").Value = "0.0₃0842"
$ws.Range("E33").Value = "  +3.01%  "
$ws.Range("E34").Value = "  -0.56%  "
$ws.Range("D35").Value = "'5.22"
$ws.Range("E35").Value = "  +0.54%  "
$ws.Range("D36").Value = "'169.01"
$ws.Range("E36").Value = "  +0.58%  "
$ws.Range("D37").Value = "'0.407"
$ws.Range("E37").Value = "  -0.11%  "
$ws.Range("E38").Value = "  -0.12%  "
$ws.Range("E39").Value = "  +0.63%  "
$ws.Range("E40").Value = "  +1.13%  "
$ws.Range("E41").Value = "  +0.03%  "
$ws.Range("D42").Value = "'168.01"
$ws.Range("E42").Value = "  -0.30%  "
$ws.Range("D43").Value = "'39.81"
$ws.Range("E43").Value = "  -0.43%  "
$ws.Range("D44").Value = "'3.91"
$ws.Range("E44").Value = "  +3.70%  "
$ws.Range("D45").Value = "'0.0596"
$ws.Range("E45").Value = "  +4.39%  "
$ws.Range("D46").Value = "'21.30"
$ws.Range("E46").Value = "  -4.00%  "
$ws.Range("D47").Value = "'0.626"
$ws.Range("E47").Value = "  -0.59%  "
$ws.Range("E48").Value = "  +0.82%  "
$ws.Range("E49").Value = "  +5.77%  "
$ws.Range("D50").Value = "'0.0964"
$ws.Range("E50").Value = "  +0.37%  "
$ws.Range("D51").Value = "'19.11"
$ws.Range("E51").Value = "  +1.66%  "
